$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename "Mælk" -> "Mælkebøtte" in A2
$ws.Range("A2").Value = "Mælkebøtte"

# 2) Add new row 9: Flute / 800 / 200 / 100 / 10
$ws.Range("A9").Value = "Flute"
$ws.Range("B9").Value = 800
$ws.Range("C9").Value = 200
$ws.Range("D9").Value = 100
$ws.Range("E9").Value = 10

Write-Host "done"
